$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# New data table (CLINIC, RESPONSE, COMMENTS) replacing the existing rows 2-9
# and adding two additional rows (10-11).
$data = @(
    @("Heart Failure", "Extremely Unlikely", "Doctors are patronising and made me feel bad"),
    @("Heart Failure", "Extremely Unlikely", "Felt as if i was not a priority"),
    @("A&E", "Extremely Unlikely", "Doctors are patronising and made me feel bad"),
    @("A&E", "Unlikely", "Waited too long to find a parking spot"),
    @("Bone Health", "Extremely Unlikely", "Waited too long to find a parking spot"),
    @("A&E", "Extremely Unlikely", "Waited over 5 hours"),
    @("Theatre Treatment Suite Implants", "Extremely Unlikely", "Felt as if i was not a priority"),
    @("Labour and Delivery Suite", "Extremely Unlikely", "Doctors are patronising and made me feel bad"),
    @("A&E", "Extremely Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("Labour and Delivery Suite", "Unlikely", "Clenliness isn't the best but otherwise okay")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("A12:XFD94").Select()
